$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" date column (C2:C11) from 45171 (2023-09-02) to 45172 (2023-09-03)
for ($row = 2; $row -le 11; $row++) {
    $ws.Cells.Item($row, 3).Value = 45172
}

# Update Signalarter (I2) from 5 to 6
$ws.Cells.Item(2, 9).Value = 6

# Update Alla arter (Q2) from 9 to 10
$ws.Cells.Item(2, 17).Value = 10

# Update Artnamn (R2): insert "Grovticka" after "Granbarkgnagare" and before "Mindre märgborre"
$ws.Cells.Item(2, 18).Value = "Spillkråka`r`nTallticka`r`nBjörksplintborre`r`nBlåmossa`r`nBronshjon`r`nGranbarkgnagare`r`nGrovticka`r`nMindre märgborre`r`nFläcknycklar`r`nRevlummer"
